{"js": "// Replace the date line and the 24 multiplication problems with their\n// updated values. Each old value is unique within the document, so a\n// direct search-and-replace for each pair is safe and keeps all\n// existing run/paragraph formatting intact.\nconst replacements = [\n  [\"2025-12-05 Friday\", \"2025-12-06 Saturday\"],\n  [\"581\u00d73=\", \"522\u00d73=\"],\n  [\"385\u00d72=\", \"319\u00d79=\"],\n  [\"961\u00d79=\", \"514\u00d73=\"],\n  [\"366\u00d78=\", \"967\u00d75=\"],\n  [\"768\u00d79=\", \"766\u00d72=\"],\n  [\"927\u00d76=\", \"421\u00d79=\"],\n  [\"771\u00d74=\", \"541\u00d73=\"],\n  [\"468\u00d76=\", \"242\u00d76=\"],\n  [\"963\u00d72=\", \"466\u00d74=\"],\n  [\"396\u00d76=\", \"105\u00d74=\"],\n  [\"841\u00d78=\", \"823\u00d72=\"],\n  [\"892\u00d76=\", \"537\u00d72=\"],\n  [\"518\u00d79=\", \"782\u00d75=\"],\n  [\"395\u00d75=\", \"706\u00d77=\"],\n  [\"960\u00d78=\", \"910\u00d78=\"],\n  [\"508\u00d77=\", \"632\u00d73=\"],\n  [\"632\u00d79=\", \"860\u00d77=\"],\n  [\"422\u00d77=\", \"476\u00d72=\"],\n  [\"203\u00d75=\", \"939\u00d79=\"],\n  [\"843\u00d73=\", \"484\u00d79=\"],\n  [\"646\u00d73=\", \"142\u00d78=\"],\n  [\"594\u00d72=\", \"806\u00d79=\"],\n  [\"686\u00d73=\", \"970\u00d75=\"],\n  [\"648\u00d72=\", \"500\u00d74=\"],\n  [\"571\u00d72=\", \"624\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 24 multiplication problems with their\n# updated values. Each old value is unique within the document, so a\n# direct Find/Replace for each pair is safe and preserves all existing\n# run/paragraph formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-05 Friday\", \"2025-12-06 Saturday\"),\n    @(\"581\u00d73=\", \"522\u00d73=\"),\n    @(\"385\u00d72=\", \"319\u00d79=\"),\n    @(\"961\u00d79=\", \"514\u00d73=\"),\n    @(\"366\u00d78=\", \"967\u00d75=\"),\n    @(\"768\u00d79=\", \"766\u00d72=\"),\n    @(\"927\u00d76=\", \"421\u00d79=\"),\n    @(\"771\u00d74=\", \"541\u00d73=\"),\n    @(\"468\u00d76=\", \"242\u00d76=\"),\n    @(\"963\u00d72=\", \"466\u00d74=\"),\n    @(\"396\u00d76=\", \"105\u00d74=\"),\n    @(\"841\u00d78=\", \"823\u00d72=\"),\n    @(\"892\u00d76=\", \"537\u00d72=\"),\n    @(\"518\u00d79=\", \"782\u00d75=\"),\n    @(\"395\u00d75=\", \"706\u00d77=\"),\n    @(\"960\u00d78=\", \"910\u00d78=\"),\n    @(\"508\u00d77=\", \"632\u00d73=\"),\n    @(\"632\u00d79=\", \"860\u00d77=\"),\n    @(\"422\u00d77=\", \"476\u00d72=\"),\n    @(\"203\u00d75=\", \"939\u00d79=\"),\n    @(\"843\u00d73=\", \"484\u00d79=\"),\n    @(\"646\u00d73=\", \"142\u00d78=\"),\n    @(\"594\u00d72=\", \"806\u00d79=\"),\n    @(\"686\u00d73=\", \"970\u00d75=\"),\n    @(\"648\u00d72=\", \"500\u00d74=\"),\n    @(\"571\u00d72=\", \"624\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
